$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 146.620486
$ws.Cells.Item(2, 8).Value = 439.861458
$ws.Cells.Item(2, 9).Value = 0.3983053592962091
$ws.Cells.Item(2, 10).Value = 0.3983053592962091
$ws.Cells.Item(2, 13).Value = 8.820647333333334
$ws.Cells.Item(2, 14).Value = 26.461942
$ws.Cells.Item(2, 15).Value = 0.06415146660411865
$ws.Cells.Item(2, 16).Value = 0.06415146660411865
$ws.Cells.Item(2, 17).Value = 1293.287598847937
$ws.Cells.Item(2, 18).Value = 11639.58838963144
$ws.Cells.Item(2, 19).Value = 0.02555187295513224
$ws.Cells.Item(2, 20).Value = 0.02555187295513224

$ws.Cells.Item(3, 7).Value = 146.620486
$ws.Cells.Item(3, 8).Value = 439.861458
$ws.Cells.Item(3, 9).Value = 0.3983053592962091
$ws.Cells.Item(3, 10).Value = 0.3983053592962091
$ws.Cells.Item(3, 15).Value = 0.3979101621202897
$ws.Cells.Item(3, 16).Value = 0.3979101621202898
$ws.Cells.Item(3, 17).Value = 8021.83185150601
$ws.Cells.Item(3, 18).Value = 72196.48666355408
$ws.Cells.Item(3, 19).Value = 0.1584897500909348
$ws.Cells.Item(3, 20).Value = 0.1584897500909348

$ws.Cells.Item(4, 7).Value = 146.620486
$ws.Cells.Item(4, 8).Value = 439.861458
$ws.Cells.Item(4, 9).Value = 0.3983053592962091
$ws.Cells.Item(4, 10).Value = 0.3983053592962091
$ws.Cells.Item(4, 13).Value = 21.90816333333333
$ws.Cells.Item(4, 14).Value = 65.72449
$ws.Cells.Item(4, 15).Value = 0.1593353362087987
$ws.Cells.Item(4, 16).Value = 0.1593353362087987
$ws.Cells.Item(4, 17).Value = 3212.185555300714
$ws.Cells.Item(4, 18).Value = 28909.66999770642
$ws.Cells.Item(4, 19).Value = 0.06346411833722783
$ws.Cells.Item(4, 20).Value = 0.06346411833722783

$ws.Cells.Item(5, 7).Value = 146.620486
$ws.Cells.Item(5, 8).Value = 439.861458
$ws.Cells.Item(5, 9).Value = 0.3983053592962091
$ws.Cells.Item(5, 10).Value = 0.3983053592962091
$ws.Cells.Item(5, 13).Value = 52.056859
$ws.Cells.Item(5, 14).Value = 156.170577
$ws.Cells.Item(5, 15).Value = 0.3786030350667928
$ws.Cells.Item(5, 16).Value = 0.3786030350667929
$ws.Cells.Item(5, 17).Value = 7632.601966213473
$ws.Cells.Item(5, 18).Value = 68693.41769592125
$ws.Cells.Item(5, 19).Value = 0.1507996179129142
$ws.Cells.Item(5, 20).Value = 0.1507996179129142

$ws.Cells.Item(6, 9).Value = 0.534552907532962
$ws.Cells.Item(6, 10).Value = 0.5345529075329621
$ws.Cells.Item(6, 13).Value = 8.820647333333334
$ws.Cells.Item(6, 14).Value = 26.461942
$ws.Cells.Item(6, 15).Value = 0.06415146660411865
$ws.Cells.Item(6, 16).Value = 0.06415146660411865
$ws.Cells.Item(6, 17).Value = 1735.68000054542
$ws.Cells.Item(6, 18).Value = 15621.12000490878
$ws.Cells.Item(6, 19).Value = 0.03429235299573534
$ws.Cells.Item(6, 20).Value = 0.03429235299573535

$ws.Cells.Item(7, 9).Value = 0.534552907532962
$ws.Cells.Item(7, 10).Value = 0.5345529075329621
$ws.Cells.Item(7, 15).Value = 0.3979101621202897
$ws.Cells.Item(7, 16).Value = 0.3979101621202898
$ws.Cells.Item(7, 19).Value = 0.2127040340983132
$ws.Cells.Item(7, 20).Value = 0.2127040340983132

$ws.Cells.Item(8, 9).Value = 0.534552907532962
$ws.Cells.Item(8, 10).Value = 0.5345529075329621
$ws.Cells.Item(8, 13).Value = 21.90816333333333
$ws.Cells.Item(8, 14).Value = 65.72449
$ws.Cells.Item(8, 15).Value = 0.1593353362087987
$ws.Cells.Item(8, 16).Value = 0.1593353362087987
$ws.Cells.Item(8, 17).Value = 4310.971690552699
$ws.Cells.Item(8, 18).Value = 38798.74521497429
$ws.Cells.Item(8, 19).Value = 0.08517316724315536
$ws.Cells.Item(8, 20).Value = 0.08517316724315539

$ws.Cells.Item(9, 9).Value = 0.534552907532962
$ws.Cells.Item(9, 10).Value = 0.5345529075329621
$ws.Cells.Item(9, 13).Value = 52.056859
$ws.Cells.Item(9, 14).Value = 156.170577
$ws.Cells.Item(9, 15).Value = 0.3786030350667928
$ws.Cells.Item(9, 16).Value = 0.3786030350667929
$ws.Cells.Item(9, 17).Value = 10243.47144183668
$ws.Cells.Item(9, 18).Value = 92191.24297653011
$ws.Cells.Item(9, 19).Value = 0.2023833531957581
$ws.Cells.Item(9, 20).Value = 0.2023833531957581

$ws.Cells.Item(10, 7).Value = 24.174389
$ws.Cells.Item(10, 8).Value = 72.523167
$ws.Cells.Item(10, 9).Value = 0.0656715098899026
$ws.Cells.Item(10, 10).Value = 0.0656715098899026
$ws.Cells.Item(10, 13).Value = 8.820647333333334
$ws.Cells.Item(10, 14).Value = 26.461942
$ws.Cells.Item(10, 15).Value = 0.06415146660411865
$ws.Cells.Item(10, 16).Value = 0.06415146660411865
$ws.Cells.Item(10, 17).Value = 213.2337598678127
$ws.Cells.Item(10, 18).Value = 1919.103838810314
$ws.Cells.Item(10, 19).Value = 0.004212923673544134
$ws.Cells.Item(10, 20).Value = 0.004212923673544134

$ws.Cells.Item(11, 7).Value = 24.174389
$ws.Cells.Item(11, 8).Value = 72.523167
$ws.Cells.Item(11, 9).Value = 0.0656715098899026
$ws.Cells.Item(11, 10).Value = 0.0656715098899026
$ws.Cells.Item(11, 15).Value = 0.3979101621202897
$ws.Cells.Item(11, 16).Value = 0.3979101621202898
$ws.Cells.Item(11, 17).Value = 1322.617929877115
$ws.Cells.Item(11, 18).Value = 11903.56136889403
$ws.Cells.Item(11, 19).Value = 0.02613136114697535
$ws.Cells.Item(11, 20).Value = 0.02613136114697536

$ws.Cells.Item(12, 7).Value = 24.174389
$ws.Cells.Item(12, 8).Value = 72.523167
$ws.Cells.Item(12, 9).Value = 0.0656715098899026
$ws.Cells.Item(12, 10).Value = 0.0656715098899026
$ws.Cells.Item(12, 13).Value = 21.90816333333333
$ws.Cells.Item(12, 14).Value = 65.72449
$ws.Cells.Item(12, 15).Value = 0.1593353362087987
$ws.Cells.Item(12, 16).Value = 0.1593353362087987
$ws.Cells.Item(12, 17).Value = 529.6164626955367
$ws.Cells.Item(12, 18).Value = 4766.548164259831
$ws.Cells.Item(12, 19).Value = 0.01046379210764708
$ws.Cells.Item(12, 20).Value = 0.01046379210764708

$ws.Cells.Item(13, 7).Value = 24.174389
$ws.Cells.Item(13, 8).Value = 72.523167
$ws.Cells.Item(13, 9).Value = 0.0656715098899026
$ws.Cells.Item(13, 10).Value = 0.0656715098899026
$ws.Cells.Item(13, 13).Value = 52.056859
$ws.Cells.Item(13, 14).Value = 156.170577
$ws.Cells.Item(13, 15).Value = 0.3786030350667928
$ws.Cells.Item(13, 16).Value = 0.3786030350667929
$ws.Cells.Item(13, 17).Value = 1258.442759584151
$ws.Cells.Item(13, 18).Value = 11325.98483625736
$ws.Cells.Item(13, 19).Value = 0.02486343296173603
$ws.Cells.Item(13, 20).Value = 0.02486343296173603

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.541205
$ws.Cells.Item(14, 8).Value = 1.623615
$ws.Cells.Item(14, 9).Value = 0.001470223280926138
$ws.Cells.Item(14, 10).Value = 0.001470223280926138
$ws.Cells.Item(14, 13).Value = 8.820647333333334
$ws.Cells.Item(14, 14).Value = 26.461942
$ws.Cells.Item(14, 15).Value = 0.06415146660411865
$ws.Cells.Item(14, 16).Value = 0.06415146660411865
$ws.Cells.Item(14, 17).Value = 4.773778440036667
$ws.Cells.Item(14, 18).Value = 42.96400596033
$ws.Cells.Item(14, 19).Value = 0.00009431697970693088
$ws.Cells.Item(14, 20).Value = 0.00009431697970693088

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.541205
$ws.Cells.Item(15, 8).Value = 1.623615
$ws.Cells.Item(15, 9).Value = 0.001470223280926138
$ws.Cells.Item(15, 10).Value = 0.001470223280926138
$ws.Cells.Item(15, 15).Value = 0.3979101621202897
$ws.Cells.Item(15, 16).Value = 0.3979101621202898
$ws.Cells.Item(15, 17).Value = 29.610156299675
$ws.Cells.Item(15, 18).Value = 266.491406697075
$ws.Cells.Item(15, 19).Value = 0.0005850167840663438
$ws.Cells.Item(15, 20).Value = 0.0005850167840663439

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.541205
$ws.Cells.Item(16, 8).Value = 1.623615
$ws.Cells.Item(16, 9).Value = 0.001470223280926138
$ws.Cells.Item(16, 10).Value = 0.001470223280926138
$ws.Cells.Item(16, 13).Value = 21.90816333333333
$ws.Cells.Item(16, 14).Value = 65.72449
$ws.Cells.Item(16, 15).Value = 0.1593353362087987
$ws.Cells.Item(16, 16).Value = 0.1593353362087987
$ws.Cells.Item(16, 17).Value = 11.85680753681667
$ws.Cells.Item(16, 18).Value = 106.71126783135
$ws.Cells.Item(16, 19).Value = 0.0002342585207683692
$ws.Cells.Item(16, 20).Value = 0.0002342585207683692

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.541205
$ws.Cells.Item(17, 8).Value = 1.623615
$ws.Cells.Item(17, 9).Value = 0.001470223280926138
$ws.Cells.Item(17, 10).Value = 0.001470223280926138
$ws.Cells.Item(17, 13).Value = 52.056859
$ws.Cells.Item(17, 14).Value = 156.170577
$ws.Cells.Item(17, 15).Value = 0.3786030350667928
$ws.Cells.Item(17, 16).Value = 0.3786030350667929
$ws.Cells.Item(17, 17).Value = 28.173432375095
$ws.Cells.Item(17, 18).Value = 253.560891375855
$ws.Cells.Item(17, 19).Value = 0.0005566309963844938
$ws.Cells.Item(17, 20).Value = 0.0005566309963844939
